# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" with the same
#    layout as the other quarterly holdings sheets: 基金代码/基金名称/基金规模/
#    股票总仓位/仓位占比/持有市值(亿元)/仓位排名.
# 2. Prepend a "2022-Q1" row to the "总计" (totals) summary sheet, pushing the
#    existing 2021-Q4 / 2021-Q3 rows down one row.

$xlPasteFormats = -4122
$xlCenter = -4108
$xlTop = -4160
$xlContinuous = 1

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Build the new "2022-Q1" sheet (cloned layout/formatting from 2021-Q4)
# ---------------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $q4Sheet)
$newSheet.Name = "2022-Q1"

# Pull over the header/index-column formatting (bold, centered, bordered)
# from the 2021-Q4 sheet so the new sheet matches the existing visual style.
$q4Sheet.Range("A1:H17").Copy()
$newSheet.Range("A1:H17").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $newSheet.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

$rows = @(
    @(0, "150103", "银河银泰混合", "15.32", "76.40", "4.21", "0.6450", 7),
    @(1, "519670", "银河行业混合", "9.93", "80.84", "3.50", "0.3476", 6),
    @(2, "519679", "银河主题混合", "6.89", "90.38", "4.84", "0.3335", 5),
    @(3, "005823", "泰康颐享混合A", "14.39", "20.19", "1.98", "0.2849", 3),
    @(4, "151001", "银河稳健混合", "8.69", "72.74", "2.81", "0.2442", 10),
    @(5, "166011", "中欧盛世成长混合 (LOF) -A", "5.80", "85.98", "3.28", "0.1902", 8),
    @(6, "001888", "中欧盛世成长混合 (LOF) -E", "5.80", "85.98", "3.28", "0.1902", 8),
    @(7, "001306", "中欧永裕混合A", "4.48", "86.33", "3.29", "0.1474", 8),
    @(8, "007203", "银河新动能混合", "3.10", "89.54", "4.43", "0.1373", 7),
    @(9, "009490", "泰康科技创新一年定期开放混合", "2.61", "79.69", "4.87", "0.1271", 7),
    @(10, "519668", "银河竞争优势成长混合", "2.12", "86.14", "4.69", "0.0994", 7),
    @(11, "519642", "银河大国智造主题灵活配置混合", "2.31", "90.40", "4.25", "0.0982", 9),
    @(12, "005824", "泰康颐享混合C", "2.82", "20.19", "1.98", "0.0558", 3),
    @(13, "008709", "银河龙头精选股票", "0.97", "81.72", "3.81", "0.0370", 8),
    @(14, "004233", "中欧盛世成长混合 (LOF) -C", "0.44", "85.98", "3.28", "0.0144", 8),
    @(15, "001307", "中欧永裕混合C", "0.35", "86.33", "3.29", "0.0115", 8)
)

$r = 2
foreach ($row in $rows) {
    $newSheet.Cells.Item($r, 1).Value = $row[0]

    $c2 = $newSheet.Cells.Item($r, 2)
    $c2.NumberFormat = "@"
    $c2.Value = $row[1]

    $newSheet.Cells.Item($r, 3).Value = $row[2]

    $c4 = $newSheet.Cells.Item($r, 4)
    $c4.NumberFormat = "@"
    $c4.Value = $row[3]

    $c5 = $newSheet.Cells.Item($r, 5)
    $c5.NumberFormat = "@"
    $c5.Value = $row[4]

    $c6 = $newSheet.Cells.Item($r, 6)
    $c6.NumberFormat = "@"
    $c6.Value = $row[5]

    $c7 = $newSheet.Cells.Item($r, 7)
    $c7.NumberFormat = "@"
    $c7.Value = $row[6]

    $newSheet.Cells.Item($r, 8).Value = $row[7]

    $r = $r + 1
}

$newSheet.Range("A1").Select()

# ---------------------------------------------------------------------------
# 2) Prepend the 2022-Q1 row to the "总计" summary sheet
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")

# Grab the existing index-column (A) formatting before the insert shifts rows
# down, then re-apply it to the freshly inserted row so every "A" cell keeps
# the same bold/centered/bordered look.
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("A2:D2").ClearFormats()
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 16
$totalSheet.Cells.Item(2, 4).Value = 2.96

# Renumber the index column (A) for the rows pushed down.
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2

# Restore the originally active sheet/tab.
$wb.Worksheets.Item("2021-Q3").Activate()
